$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.397.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.05%  "

$ws.Range("D3").Value = "'1.848.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "

$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Value = "'240.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.14%  "

$ws.Range("D6").Value = "'0.6290"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.27%  "

$ws.Range("D8").Value = "'0.07501"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.66%  "

$ws.Range("D9").Value = "'0.2893"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.95%  "

$ws.Range("D10").Value = "'24.40"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.52%  "

$ws.Range("D11").Value = "'0.07748"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.01%  "

$ws.Range("D12").Value = "'1.848.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.80%  "

$ws.Range("D13").Value = "'5.013"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.43%  "

$ws.Range("D14").Value = "'0.6803"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").Value = "'0.00001034"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.82%  "

$ws.Range("D16").Value = "'82.88"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.83%  "

$ws.Range("D17").Value = "'2.110.52"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.01%  "

$ws.Range("D18").Value = "'6.112"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.28%  "

$ws.Range("D19").Value = "'29.422.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.03%  "

$ws.Range("D20").Value = "'229.24"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "

$ws.Range("D21").Value = "'12.31"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.70%  "

$ws.Range("E22").Value = "  +0.19%  "

$ws.Range("D23").Value = "'7.434"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.28%  "

$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").Value = "'159.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.00%  "

$ws.Range("D26").Value = "'0.1386"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.50%  "

$ws.Range("D27").Value = "'8.413"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "

$ws.Range("D28").Value = "'17.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.58%  "

$ws.Range("D29").Value = "'1.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.45%  "

$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("E31").Value = "  +0.89%  "

$ws.Range("D32").Value = "'4.120"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.12%  "

$ws.Range("D33").Value = "'4.051"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.11%  "

$ws.Range("D34").Value = "'1.823"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("E35").Value = "  -1.28%  "

$ws.Range("D36").Value = "'0.6945"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.16%  "

$ws.Range("D37").Value = "'2.590"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "

$ws.Range("D38").Value = "'2.844"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.20%  "

$ws.Range("D39").Value = "'1.250.73"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.00%  "

$ws.Range("D40").Value = "'0.01819"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.64%  "

$ws.Range("D41").Value = "'6.498"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.79%  "

$ws.Range("D42").Value = "'0.9056"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.20%  "

$ws.Range("D43").Value = "'1.002"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "

$ws.Range("D44").Value = "'2.011.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.27%  "

$ws.Range("D45").Value = "'101.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.40%  "

$ws.Range("D46").Value = "'65.83"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.39%  "

$ws.Range("D47").Value = "'7.089"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.39%  "

$ws.Range("E48").Value = "  -0.02%  "

$ws.Range("E49").Value = "  -4.72%  "

$ws.Range("D50").Value = "'8.957"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.59%  "

$ws.Range("E51").Value = "  -1.99%  "

